# Applies the Tue Nov 14 15:40:58 UTC 2023 "cryptos" data refresh
# (coin prices / 1h volume changes, plus a few rank swaps) to sheet1.
# D-column prices are prefixed with a literal apostrophe so Excel keeps
# them as text (matching the original sheet) instead of auto-coercing
# number-like strings (e.g. "1.00", "0.660") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'36.122.14"
$ws.Cells.Item(2, 5).Value = "  -1.81%  "

$ws.Cells.Item(3, 4).Value = "'2.032.71"
$ws.Cells.Item(3, 5).Value = "  -2.83%  "

$ws.Cells.Item(4, 5).Value = "  +0.17%  "

$ws.Cells.Item(5, 4).Value = "'243.76"
$ws.Cells.Item(5, 5).Value = "  -1.00%  "

$ws.Cells.Item(6, 4).Value = "'0.660"
$ws.Cells.Item(6, 5).Value = "  +1.28%  "

$ws.Cells.Item(7, 5).Value = "  +0.13%  "

$ws.Cells.Item(8, 4).Value = "'55.90"
$ws.Cells.Item(8, 5).Value = "  +0.40%  "

$ws.Cells.Item(9, 4).Value = "'62.38"
$ws.Cells.Item(9, 5).Value = "  +4.68%  "

$ws.Cells.Item(10, 4).Value = "'0.363"
$ws.Cells.Item(10, 5).Value = "  -1.36%  "

$ws.Cells.Item(11, 4).Value = "'0.0739"
$ws.Cells.Item(11, 5).Value = "  -3.90%  "

$ws.Cells.Item(12, 5).Value = "  -3.28%  "

$ws.Cells.Item(13, 4).Value = "'0.894"
$ws.Cells.Item(13, 5).Value = "  +0.97%  "

$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).Value = "'14.05"
$ws.Cells.Item(14, 5).Value = "  -5.96%  "

$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "'2.333.74"
$ws.Cells.Item(15, 5).Value = "  -2.61%  "

$ws.Cells.Item(16, 5).Value = "  -4.19%  "

$ws.Cells.Item(17, 4).Value = "'2.027.55"
$ws.Cells.Item(17, 5).Value = "  -3.04%  "

$ws.Cells.Item(18, 4).Value = "'36.126.63"
$ws.Cells.Item(18, 5).Value = "  -1.79%  "

$ws.Cells.Item(19, 4).Value = "'17.25"

$ws.Cells.Item(20, 4).Value = "'71.08"
$ws.Cells.Item(20, 5).Value = "  -2.90%  "

$ws.Cells.Item(21, 4).Value = "'0.0₃0849"
$ws.Cells.Item(21, 5).Value = "  -3.47%  "

$ws.Cells.Item(22, 4).Value = "'236.31"
$ws.Cells.Item(22, 5).Value = "  -0.14%  "

$ws.Cells.Item(23, 4).Value = "'5.13"
$ws.Cells.Item(23, 5).Value = "  -6.42%  "

$ws.Cells.Item(24, 5).Value = "  -0.29%  "

$ws.Cells.Item(25, 5).Value = "  -2.91%  "

$ws.Cells.Item(26, 5).Value = "  +2.65%  "

$ws.Cells.Item(27, 4).Value = "'9.15"
$ws.Cells.Item(27, 5).Value = "  -8.09%  "

$ws.Cells.Item(28, 4).Value = "'163.18"
$ws.Cells.Item(28, 5).Value = "  -3.04%  "

$ws.Cells.Item(29, 4).Value = "'19.82"
$ws.Cells.Item(29, 5).Value = "  -5.63%  "

$ws.Cells.Item(30, 5).Value = "  -2.93%  "

$ws.Cells.Item(31, 4).Value = "'1.19"
$ws.Cells.Item(31, 5).Value = "  -2.02%  "

$ws.Cells.Item(32, 5).Value = "  -7.65%  "

$ws.Cells.Item(33, 4).Value = "'0.0595"
$ws.Cells.Item(33, 5).Value = "  -2.56%  "

$ws.Cells.Item(34, 5).Value = "  -7.50%  "

$ws.Cells.Item(35, 2).Value = "BinanceUSD"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 5).Value = "  +0.06%  "

$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).Value = "'0.0863"
$ws.Cells.Item(36, 5).Value = "  +3.02%  "

$ws.Cells.Item(37, 5).Value = "  -0.80%  "

$ws.Cells.Item(38, 4).Value = "'2.18"
$ws.Cells.Item(38, 5).Value = "  -9.82%  "

$ws.Cells.Item(39, 4).Value = "'5.00"
$ws.Cells.Item(39, 5).Value = "  +1.30%  "

$ws.Cells.Item(40, 4).Value = "'1.20"
$ws.Cells.Item(40, 5).Value = "  -6.58%  "

$ws.Cells.Item(41, 5).Value = "  -1.47%  "

$ws.Cells.Item(42, 5).Value = "  -3.54%  "

$ws.Cells.Item(43, 5).Value = "  -6.75%  "

$ws.Cells.Item(44, 4).Value = "'92.58"
$ws.Cells.Item(44, 5).Value = "  -4.42%  "

$ws.Cells.Item(45, 4).Value = "'0.0895"
$ws.Cells.Item(45, 5).Value = "  -6.45%  "

$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(46, 4).Value = "'1.353.72"
$ws.Cells.Item(46, 5).Value = "  +0.90%  "

$ws.Cells.Item(47, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(47, 4).Value = "'15.55"
$ws.Cells.Item(47, 5).Value = "  -5.03%  "

$ws.Cells.Item(48, 4).Value = "'7.34"
$ws.Cells.Item(48, 5).Value = "  +3.95%  "

$ws.Cells.Item(49, 4).Value = "'2.92"
$ws.Cells.Item(49, 5).Value = "  +1.36%  "

$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(50, 4).Value = "'2.23"
$ws.Cells.Item(50, 5).Value = "  -8.42%  "

$ws.Cells.Item(51, 4).Value = "'45.29"
$ws.Cells.Item(51, 5).Value = "  -1.15%  "
